# Updated symbol list refresh: update Price (D) and Volume(1h) (E) columns
# for the cryptos sheet. Values are written with a leading apostrophe so
# Excel stores them as literal text (matching the original inlineStr cells)
# instead of auto-converting numeric-looking / percent-looking strings into
# Number/Percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.61"
$ws.Range("E2").Value = "'-1.02%"

$ws.Range("D3").Value = "'31.46"
$ws.Range("E3").Value = "'-1.09%"

$ws.Range("D4").Value = "'5.084"
$ws.Range("E4").Value = "'-1.50%"

$ws.Range("D5").Value = "'0.07853"
$ws.Range("E5").Value = "'-0.61%"

$ws.Range("D6").Value = "'2.267"
$ws.Range("E6").Value = "'-3.06%"

$ws.Range("D7").Value = "'7.813"
$ws.Range("E7").Value = "'-1.69%"

$ws.Range("D8").Value = "'3.848"
$ws.Range("E8").Value = "'-0.58%"

$ws.Range("D9").Value = "'0.9199"
$ws.Range("E9").Value = "'1.44%"

$ws.Range("D10").Value = "'0.1748"
$ws.Range("E10").Value = "'0.71%"

$ws.Range("D11").Value = "'0.07605"
$ws.Range("E11").Value = "'3.14%"

$ws.Range("D12").Value = "'0.09186"
$ws.Range("E12").Value = "'12.53%"

$ws.Range("D13").Value = "'0.03000"
$ws.Range("E13").Value = "'-3.41%"

$ws.Range("E14").Value = "'0.77%"

$ws.Range("D15").Value = "'0.001505"
$ws.Range("E15").Value = "'-1.21%"

$ws.Range("D16").Value = "'0.005965"
$ws.Range("E16").Value = "'-1.58%"

$ws.Range("E17").Value = "'-0.69%"

$ws.Range("E18").Value = "'0.09%"

$ws.Range("D20").Value = "'0.1318"
$ws.Range("E20").Value = "'-1.49%"

$ws.Range("D21").Value = "'3.915"
$ws.Range("E21").Value = "'-16.43%"

$ws.Range("E22").Value = "'3.87%"

$ws.Range("D23").Value = "'0.04617"
$ws.Range("E23").Value = "'-0.97%"

$ws.Range("D24").Value = "'0.001253"
$ws.Range("E24").Value = "'-0.56%"

$ws.Range("D25").Value = "'0.004469"
$ws.Range("E25").Value = "'-1.04%"

$ws.Range("D26").Value = "'0.0001249"
$ws.Range("E26").Value = "'-7.42%"

$ws.Range("D27").Value = "'0.0003398"
$ws.Range("E27").Value = "'23.92%"

$ws.Range("D39").Value = "'0.01733"
$ws.Range("E39").Value = "'-10.55%"

$ws.Range("D40").Value = "'0.04620"
$ws.Range("E40").Value = "'1.07%"

$ws.Range("D41").Value = "'0.007058"
$ws.Range("E41").Value = "'-3.35%"

$ws.Range("E42").Value = "'-0.22%"

$ws.Range("D43").Value = "'0.002188"
$ws.Range("E43").Value = "'-2.70%"

$ws.Range("D44").Value = "'0.009737"
$ws.Range("E44").Value = "'-9.36%"

$ws.Range("D45").Value = "'0.00006282"
$ws.Range("E45").Value = "'-3.07%"

$ws.Range("E46").Value = "'0.06%"

$ws.Range("D48").Value = "'0.7461"
$ws.Range("E48").Value = "'-9.07%"

$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.06%"

$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.06%"
